$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 8 de Octubre de 2020 a las 01:12"

# Row 4: Estados Unidos (updated counts)
$ws.Range("B4").Value = 7767598
$ws.Range("C4").Value = 39968
$ws.Range("D4").Value = 4970847
$ws.Range("E4").Value = 2580150
$ws.Range("G4").Value = 749
$ws.Range("H4").Value = 216601

# Row 8: España -> Colombia
$ws.Range("A8").Value = "Colombia"
$ws.Range("B8").Value = 877683
$ws.Range("C8").Value = 7875
$ws.Range("D8").Value = 773973
$ws.Range("E8").Value = 76530
$ws.Range("G8").Value = 163
$ws.Range("H8").Value = 27180

# Row 9: Colombia -> España
$ws.Range("A9").Value = "España"
$ws.Range("B9").Value = 872276
$ws.Range("C9").Value = 6645
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("G9").Value = 76
$ws.Range("H9").Value = 32562

# Row 10: Peru -> Argentina
$ws.Range("A10").Value = "Argentina"
$ws.Range("B10").Value = 840915
$ws.Range("C10").Value = 16447
$ws.Range("D10").Value = 670725
$ws.Range("E10").Value = 147964
$ws.Range("G10").Value = 399
$ws.Range("H10").Value = 22226

# Row 11: Argentina -> Peru
$ws.Range("A11").Value = "Peru"
$ws.Range("B11").Value = 832929
$ws.Range("D11").Value = 718065
$ws.Range("E11").Value = 81950
$ws.Range("H11").Value = 32914

# Row 26: Alemania (updated counts)
$ws.Range("B26").Value = 311113
$ws.Range("C26").Value = 3994
$ws.Range("E26").Value = 33761

# Row 37: Panama (updated counts)
$ws.Range("B37").Value = 117300
$ws.Range("C37").Value = 698
$ws.Range("D37").Value = 93610
$ws.Range("E37").Value = 21242
$ws.Range("G37").Value = 8
$ws.Range("H37").Value = 2448

# Row 42: Egipto (updated counts)
$ws.Range("B42").Value = 104035
$ws.Range("C42").Value = 133
$ws.Range("D42").Value = 97492
$ws.Range("E42").Value = 533
$ws.Range("G42").Value = 9
$ws.Range("H42").Value = 6010

# Row 58: Nigeria (updated counts)
$ws.Range("B58").Value = 59738
$ws.Range("C58").Value = 155
$ws.Range("D58").Value = 51403
$ws.Range("E58").Value = 7222

# Row 72: Kenia (updated counts)
$ws.Range("D72").Value = 31659
$ws.Range("E72").Value = 7500

# Row 99: Sudan (updated counts)
$ws.Range("B99").Value = 13668
$ws.Range("C99").Value = 15
$ws.Range("E99").Value = 6068

# Row 114: Zimbabue (updated counts)
$ws.Range("B114").Value = 7919
$ws.Range("C114").Value = 4
$ws.Range("D114").Value = 6441
$ws.Range("E114").Value = 1249

# Row 115: Mauritania (updated counts)
$ws.Range("B115").Value = 7535
$ws.Range("C115").Value = 6
$ws.Range("D115").Value = 7212
$ws.Range("E115").Value = 161

# Row 131: Ruanda -> Trinidad yTobago
$ws.Range("A131").Value = "Trinidad yTobago"
$ws.Range("B131").Value = 4887
$ws.Range("C131").Value = 41
$ws.Range("D131").Value = 3010
$ws.Range("E131").Value = 1793
$ws.Range("G131").Value = 1
$ws.Range("H131").Value = 84

# Row 132: Trinidad yTobago -> Ruanda
$ws.Range("A132").Value = "Ruanda"
$ws.Range("B132").Value = 4883
$ws.Range("C132").Value = 10
$ws.Range("D132").Value = 3408
$ws.Range("E132").Value = 1446
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 29

# Row 157: Burkina Faso (updated counts)
$ws.Range("B157").Value = 2222
$ws.Range("C157").Value = 25
$ws.Range("D157").Value = 1478
$ws.Range("E157").Value = 685

# Row 158: Uruguay (updated counts)
$ws.Range("B158").Value = 2206
$ws.Range("C158").Value = 29
$ws.Range("D158").Value = 1890
$ws.Range("E158").Value = 267

# Row 207: Nueva Caledonia -> Santa Lucia
$ws.Range("A207").Value = "Santa Lucia"

# Row 208: Santa Lucia -> Nueva Caledonia
$ws.Range("A208").Value = "Nueva Caledonia"

# Row 215: Montserrat -> Islas Malvinas
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0

# Row 216: Islas Malvinas -> Montserrat
$ws.Range("A216").Value = "Montserrat"
$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1
